$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 for the new fixture (Bukovyna Chernivtsi - Podillya Khmelnytskyi),
# shifting existing rows 3-21 down to rows 4-22.
$ws.Rows.Item(3).Insert()

# Row 2
$ws.Range("A2").Value = 'Sun Oct 12'
$ws.Range("B2").Value = 'Croatia ✓ - Gibraltar: 3:0'
$ws.Range("C2").Value = 3.5
$ws.Range("D2").Value = 'Croatia'
$ws.Range("E2").Value = 4.5
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '77%'
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = '✓'
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = $true

# Row 3
$ws.Range("A3").Value = 'Sun Oct 12'
$ws.Range("B3").Value = 'Bukovyna Chernivtsi ✓ - Podillya Khmelnytskyi: 2:0'
$ws.Range("C3").Value = 1.62
$ws.Range("D3").Value = 'Bukovyna Chernivtsi'
$ws.Range("E3").Value = 2.5
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = '76%'
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = '✓'
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = $true

# Row 4
$ws.Range("A4").Value = 'Sun Oct 12'
$ws.Range("B4").Value = 'Louisville City FC ✓ - Miami FC: 1:0'
$ws.Range("C4").Value = 2.02
$ws.Range("D4").Value = 'Louisville City FC'
$ws.Range("E4").Value = 3.5
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = '73%'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = '✓'
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = $true

# Row 5
$ws.Range("A5").Value = 'Sun Oct 12'
$ws.Range("B5").Value = 'Netherlands ✓ - Finland: 4:0'
$ws.Range("C5").Value = 3.03
$ws.Range("D5").Value = 'Netherlands'
$ws.Range("E5").Value = 4.5
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = '72%'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = '✓'
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = $true

# Row 6
$ws.Range("A6").Value = 'Sun Oct 12'
$ws.Range("B6").Value = 'Scotland ✓ - Belarus: 2:1'
$ws.Range("C6").Value = 1.3
$ws.Range("D6").Value = 'Scotland'
$ws.Range("E6").Value = 2.5
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = '72%'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = '✓'
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = $false

# Row 7
$ws.Range("A7").Value = 'Sun Oct 12'
$ws.Range("B7").Value = 'Romania - Austria X: 1:0'
$ws.Range("C7").Value = 1.96
$ws.Range("D7").Value = 'Austria'
$ws.Range("E7").Value = 2.5
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = '71%'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = 'X'
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = $true

# Row 8
$ws.Range("A8").Value = 'Sun Oct 12'
$ws.Range("B8").Value = 'Zambia X - Niger: 0:1'
$ws.Range("C8").Value = 0.57
$ws.Range("D8").Value = 'Zambia'
$ws.Range("E8").Value = 1.5
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = '71%'
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = 'X'
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = $true

# Row 9
$ws.Range("A9").Value = 'Sun Oct 12'
$ws.Range("B9").Value = 'Egypt ✓ - Guinea-Bissau: 1:0'
$ws.Range("C9").Value = 0.14
$ws.Range("D9").Value = 'Egypt'
$ws.Range("E9").Value = 1.5
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = '71%'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = '✓'
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = $true

# Row 10
$ws.Range("A10").Value = 'Sun Oct 12'
$ws.Range("B10").Value = 'CD Plaza Amador X - CD Árabe Unido: 0:2'
$ws.Range("C10").Value = 1.62
$ws.Range("D10").Value = 'CD Plaza Amador'
$ws.Range("E10").Value = 2.5
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = '71%'
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = 'X'
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = $true

# Row 11
$ws.Range("A11").Value = 'Sun Oct 12'
$ws.Range("B11").Value = 'Atlético Tembetary - Club Libertad Asunción X: 2:1'
$ws.Range("C11").Value = 1.34
$ws.Range("D11").Value = 'Club Libertad Asunción'
$ws.Range("E11").Value = 2.5
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = '71%'
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = 'X'
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = $false

# Row 12
$ws.Range("A12").Value = 'Sun Oct 12'
$ws.Range("B12").Value = 'Burkina Faso ✓ - Ethiopia: 3:1'
$ws.Range("C12").Value = 2.3
$ws.Range("D12").Value = 'Burkina Faso'
$ws.Range("E12").Value = 3.5
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = '70%'
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = '✓'
$ws.Range("H12").Value = 4
$ws.Range("I12").Value = $false

# Row 13
$ws.Range("A13").Value = 'Sun Oct 12'
$ws.Range("B13").Value = 'Ghana ✓ - Comoros: 1:0'
$ws.Range("C13").Value = 1.05
$ws.Range("D13").Value = 'Ghana'
$ws.Range("E13").Value = 2.5
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = '70%'
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = '✓'
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = $true

# Row 14
$ws.Range("A14").Value = 'Sun Oct 12'
$ws.Range("B14").Value = 'Mali ✓ - Madagascar: 4:1'
$ws.Range("C14").Value = 0.89
$ws.Range("D14").Value = 'Mali'
$ws.Range("E14").Value = 1.5
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = '67%'
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = '✓'
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = $false

# Row 15
$ws.Range("A15").Value = 'Sun Oct 12'
$ws.Range("B15").Value = 'Sanfrecce Hiroshima ✓ - Yokohama FC: 2:1'
$ws.Range("C15").Value = 1.5
$ws.Range("D15").Value = 'Sanfrecce Hiroshima'
$ws.Range("E15").Value = 2.5
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = '61%'
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = '✓'
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = $false

# Row 16
$ws.Range("A16").Value = 'Sun Oct 12'
$ws.Range("B16").Value = 'PSS Sleman ✓ - Kendal Tornado FC: 3:1'
$ws.Range("C16").Value = 2.12
$ws.Range("D16").Value = 'PSS Sleman'
$ws.Range("E16").Value = 3.5
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = '61%'
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = '✓'
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = $false

# Row 17
$ws.Range("A17").Value = 'Sun Oct 12'
$ws.Range("B17").Value = 'Club Deportivo Guabirá ✓ - Club Aurora: 2:1'
$ws.Range("C17").Value = 2.63
$ws.Range("D17").Value = 'Club Deportivo Guabirá'
$ws.Range("E17").Value = 3.5
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = '61%'
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = '✓'
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = $true

# Row 18
$ws.Range("A18").Value = 'Sun Oct 12'
$ws.Range("B18").Value = 'Pars Jonoubi Jam  - Fard Alborz: 0:0'
$ws.Range("C18").Value = 1.74
$ws.Range("D18").Value = 'Pars Jonoubi Jam'
$ws.Range("E18").Value = 2.5
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = '58%'
$ws.Range("F18").Style = "Normal"
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = $true

# Row 19
$ws.Range("A19").Value = 'Sun Oct 12'
$ws.Range("B19").Value = 'Antigua GFC  - Deportivo Achuapa: 02:00'
$ws.Range("C19").Value = 1.5
$ws.Range("D19").Value = 'Antigua GFC'
$ws.Range("E19").Value = 2.5
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = '55%'
$ws.Range("F19").Style = "Normal"
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = $true

# Row 20
$ws.Range("A20").Value = 'Sun Oct 12'
$ws.Range("B20").Value = 'CA Estudiantes ✓ - Club Deportivo Maipú: 1:0'
$ws.Range("C20").Value = 0.88
$ws.Range("D20").Value = 'CA Estudiantes'
$ws.Range("E20").Value = 1.5
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = '55%'
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = '✓'
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = $true

# Row 21
$ws.Range("A21").Value = 'Sun Oct 12'
$ws.Range("B21").Value = 'Jeonnam Dragons X - Ansan Greeners: 0:1'
$ws.Range("C21").Value = 1.61
$ws.Range("D21").Value = 'Jeonnam Dragons'
$ws.Range("E21").Value = 2.5
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = '55%'
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value = 'X'
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = $true

# Row 22
$ws.Range("A22").Value = 'Sun Oct 12'
$ws.Range("B22").Value = 'CSD Municipal  - Cobán Imperial: 00:00'
$ws.Range("C22").Value = 1.66
$ws.Range("D22").Value = 'CSD Municipal'
$ws.Range("E22").Value = 2.5
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = '55%'
$ws.Range("F22").Style = "Normal"
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = $true
